$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header): M1:Q1 change from inline strings ("Unnamed: N") to plain numbers 11-15 ---
$headerValues = @(11, 12, 13, 14, 15)
$headerCols = @("M", "N", "O", "P", "Q")
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $addr = $headerCols[$i] + "1"
    $ws.Range($addr).Value = $headerValues[$i]
}

# --- Column A (rows 2-121): replace Excel date-serial values (custom date/time format, style index 2) ---
# --- with plain YYYYMMDD integers, and strip the date style back to the workbook default (Normal). ---
$dateGroups = @(
    "20180208",
    "20180214",
    "20180222",
    "20180301",
    "20180308",
    "20180315",
    "20180322",
    "20180329",
    "20180405",
    "20180412",
    "20180418",
    "20180426",
    "20180503",
    "20180510",
    "20180517",
    "20180524",
    "20180531",
    "20180607",
    "20180614",
    "20180621"
)

$row = 2
foreach ($ymd in $dateGroups) {
    for ($j = 0; $j -lt 6; $j++) {
        $addr = "A" + $row
        $c = $ws.Range($addr)
        $c.Style = "Normal"
        $c.Value = [long]$ymd
        $row++
    }
}
